$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number + week-covering date range) ---
# Both strings are rich-text runs in the shared-string table where every run
# shares identical formatting, so a plain value assignment reproduces the
# same rendered text without altering the cell style.
$ws.Range("A8").Value = "Volume 32   Number  42"
$ws.Range("C9").Value = "Report Covering the Week  10/13/2025  Through  10/19/2025"

# --- Bulk numeric value updates across the weekly crime-stat table (rows 14-30) ---
$ws.Range("M14").Value = -60
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 23
$ws.Range("K15").Value = -23.333333333333
$ws.Range("L15").Value = 9.523809523809
$ws.Range("M15").Value = 35.294117647058
$ws.Range("N15").Value = -17.857142857142
$ws.Range("F16").Value = 29
$ws.Range("G16").Value = 22
$ws.Range("H16").Value = 31.818181818181
$ws.Range("I16").Value = 244
$ws.Range("J16").Value = 201
$ws.Range("K16").Value = 21.39303482587
$ws.Range("L16").Value = 6.086956521739
$ws.Range("M16").Value = 0.826446280991
$ws.Range("N16").Value = -55.555555555555
$ws.Range("C17").Value = 10
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 100
$ws.Range("F17").Value = 41
$ws.Range("G17").Value = 23
$ws.Range("H17").Value = 78.260869565217
$ws.Range("I17").Value = 444
$ws.Range("J17").Value = 385
$ws.Range("K17").Value = 15.324675324675
$ws.Range("L17").Value = 32.142857142857
$ws.Range("M17").Value = 112.44019138756
$ws.Range("N17").Value = 66.29213483146
$ws.Range("C18").Value = 6
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 100
$ws.Range("F18").Value = 21
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = 133.333333333333
$ws.Range("I18").Value = 190
$ws.Range("J18").Value = 135
$ws.Range("K18").Value = 40.74074074074
$ws.Range("L18").Value = -4.522613065326
$ws.Range("M18").Value = -37.293729372937
$ws.Range("N18").Value = -83.592400690846
$ws.Range("F19").Value = 58
$ws.Range("G19").Value = 67
$ws.Range("H19").Value = -13.432835820895
$ws.Range("I19").Value = 643
$ws.Range("J19").Value = 705
$ws.Range("K19").Value = -8.794326241134
$ws.Range("L19").Value = 19.294990723562
$ws.Range("M19").Value = 82.670454545454
$ws.Range("N19").Value = 31.762295081967
$ws.Range("C20").Value = 12
$ws.Range("D20").Value = 6
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 36
$ws.Range("H20").Value = 24.137931034482
$ws.Range("I20").Value = 429
$ws.Range("J20").Value = 366
$ws.Range("K20").Value = 17.213114754098
$ws.Range("L20").Value = 2.386634844868
$ws.Range("M20").Value = 130.645161290323
$ws.Range("N20").Value = -71.570576540755
$ws.Range("C21").Value = 48
$ws.Range("D21").Value = 38
$ws.Range("E21").Value = 26.315789473684
$ws.Range("F21").Value = 187
$ws.Range("G21").Value = 152
$ws.Range("H21").Value = 23.026315789473
$ws.Range("I21").Value = 1977
$ws.Range("J21").Value = 1825
$ws.Range("K21").Value = 8.328767123287
$ws.Range("L21").Value = 13.03602058319
$ws.Range("M21").Value = 49.886277482941
$ws.Range("N21").Value = -50.67365269461
$ws.Range("M22").Value = -23.529411764705
$ws.Range("D23").Value = 3
$ws.Range("E23").Value = -100
$ws.Range("F23").Value = 13
$ws.Range("G23").Value = 7
$ws.Range("H23").Value = 85.714285714285
$ws.Range("J23").Value = 90
$ws.Range("K23").Value = 4.444444444444
$ws.Range("L23").Value = -4.081632653061
$ws.Range("M23").Value = 77.358490566037
$ws.Range("C24").Value = 29
$ws.Range("D24").Value = 25
$ws.Range("E24").Value = 16
$ws.Range("F24").Value = 130
$ws.Range("G24").Value = 91
$ws.Range("H24").Value = 42.857142857142
$ws.Range("I24").Value = 1424
$ws.Range("J24").Value = 1029
$ws.Range("K24").Value = 38.386783284742
$ws.Range("L24").Value = 8.536585365853
$ws.Range("M24").Value = 90.120160213618
$ws.Range("C25").Value = 16
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 68
$ws.Range("G25").Value = 43
$ws.Range("H25").Value = 58.13953488372
$ws.Range("I25").Value = 584
$ws.Range("J25").Value = 399
$ws.Range("K25").Value = 46.365914786967
$ws.Range("L25").Value = 3.91459074733
$ws.Range("C26").Value = 12
$ws.Range("D26").Value = 5
$ws.Range("E26").Value = 140
$ws.Range("F26").Value = 54
$ws.Range("G26").Value = 29
$ws.Range("H26").Value = 86.206896551724
$ws.Range("I26").Value = 540
$ws.Range("J26").Value = 455
$ws.Range("K26").Value = 18.681318681318
$ws.Range("L26").Value = 25.874125874125
$ws.Range("M26").Value = 4.651162790697
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 27
$ws.Range("K27").Value = -20.588235294117
$ws.Range("L27").Value = -10
$ws.Range("C28").Value = 2
$ws.Range("E28").Value = 0
$ws.Range("G28").Value = 11
$ws.Range("H28").Value = -27.272727272727
$ws.Range("I28").Value = 59
$ws.Range("J28").Value = 55
$ws.Range("K28").Value = 7.272727272727
$ws.Range("L28").Value = 0
$ws.Range("M29").Value = -47.826086956521
$ws.Range("N29").Value = -66.666666666666
$ws.Range("M30").Value = -63.157894736842
$ws.Range("N30").Value = -77.419354838709


# --- Cells that flip between the "0" text placeholder and a real number ---
# Text("0") -> Number: copy the numeric sibling's number format first so the
# destination style matches what Excel would natively assign, then write the value.
$ws.Range("C15").NumberFormat = $ws.Range("G15").NumberFormat
$ws.Range("C15").Value = 2

$ws.Range("F15").NumberFormat = $ws.Range("G15").NumberFormat
$ws.Range("F15").Value = 2

$ws.Range("C27").NumberFormat = $ws.Range("G27").NumberFormat
$ws.Range("C27").Value = 2

$ws.Range("F27").NumberFormat = $ws.Range("G27").NumberFormat
$ws.Range("F27").Value = 2

# Number -> Text("0"): force literal text with a quote-prefix entry, then
# paste the donor cell's formatting (General/style-13) over it so the
# quote-prefix style fork collapses back onto the original text style.
$ws.Range("C23").Value = "'0"
$ws.Range("C22").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("F29").Value = "'0"
$ws.Range("G29").Copy()
$ws.Range("F29").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("F30").Value = "'0"
$ws.Range("G30").Copy()
$ws.Range("F30").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Column widths for I:J widen to match the new (wider) 28-day figures ---
$ws.Columns.Item(9).ColumnWidth = 6.719482285714285
$ws.Columns.Item(10).ColumnWidth = 6.719482285714285
